$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '30.773.51'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '2.111.47'
$ws.Range('E3').Value = '  +7.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.43'
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5292'
$ws.Range('E7').Value = '  +3.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4405'
$ws.Range('E8').Value = '  +7.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09015'
$ws.Range('E9').Value = '  +6.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.96'
$ws.Range('E10').Value = '  +8.20%  '
$ws.Range('E11').Value = '  +4.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.01'
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').Value = '2.108.90'
$ws.Range('E13').Value = '  +8.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.752'
$ws.Range('E14').Value = '  +4.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.797'
$ws.Range('E15').Value = '  +5.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.48'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06665'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.14'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  +5.66%  '
$ws.Range('D23').Value = '30.825.95'
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.25'
$ws.Range('E24').Value = '  +7.14%  '
$ws.Range('D25').Value = '2.356.62'
$ws.Range('E25').Value = '  +8.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.257'
$ws.Range('E26').Value = '  +2.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.84'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.559'
$ws.Range('E28').Value = '  +7.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.48'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.08'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  +2.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.219'
$ws.Range('E33').Value = '  +3.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.024'
$ws.Range('E34').Value = '  +6.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.546'
$ws.Range('E35').Value = '  +19.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02603'
$ws.Range('E36').Value = '  +5.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.532'
$ws.Range('E37').Value = '  +3.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06740'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.508'
$ws.Range('E39').Value = '  +7.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.67'
$ws.Range('E40').Value = '  +7.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2268'
$ws.Range('E41').Value = '  +4.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6849'
$ws.Range('E42').Value = '  +4.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.241'
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6455'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.12'
$ws.Range('E45').Value = '  +5.23%  '
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.667'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('E49').Value = '  +4.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.35'
$ws.Range('E50').Value = '  +4.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.71'
$ws.Range('E51').Value = '  -2.95%  '
